$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort rows 2-5 by year (1996, 1998, 1997, 1999) and correct the
# product_type for the non "box set" items to "supplement".

# Row 2: El Genesis (1996)
$ws.Range("A2").Value = 1996
$ws.Range("B2").Value = "創聖記エルジェネシス"
$ws.Range("C2").Value = "El Genesis"
$ws.Range("D2").Value = "Game Field"
$ws.Range("E2").Value = "el_genesis.jpg"
$ws.Range("F2").Value = "box set"

# Row 3: School Civil War: El Genesis (1998)
$ws.Range("A3").Value = 1998
$ws.Range("B3").Value = "学園戦国エルジェネシス"
$ws.Range("C3").Value = "School Civil War: El Genesis"
$ws.Range("D3").Value = "Game Field"
$ws.Range("E3").Value = "gakuen_sengoku.jpg"
$ws.Range("F3").Value = "box set"

# Row 4: Divine Era: El Genesis (1997)
$ws.Range("A4").Value = 1997
$ws.Range("B4").Value = "神世紀エルジェネシス"
$ws.Range("C4").Value = "Divine Era: El Genesis"
$ws.Range("D4").Value = "Game Field"
$ws.Range("E4").Value = "god_century.jpg"
$ws.Range("F4").Value = "box set"

# Row 5: Arena of Destiny (1999) - now a supplement, not a box set
$ws.Range("A5").Value = 1999
$ws.Range("B5").Value = "創聖記エルジェネシス2 ZWEI"
$ws.Range("C5").Value = "Arena of Destiny"
$ws.Range("D5").Value = "Game Field"
$ws.Range("E5").Value = "arena_of_destiny.jpg"
$ws.Range("F5").Value = "supplement"

# Rows 6-8 keep their content but are also supplements, not box sets
$ws.Range("F6").Value = "supplement"
$ws.Range("F7").Value = "supplement"
$ws.Range("F8").Value = "supplement"

# Update the selected/active cell to match the author's saved state
$ws.Range("F6").Select()
